$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("O2").Value = 1.5
$ws.Range("P2").Value = 2.5

# Row 4 updates
$ws.Range("G4").Value = 3.2
$ws.Range("I4").Value = 2.3
$ws.Range("L4").Value = 3.1
$ws.Range("AK4").Value = 21
$ws.Range("AX4").Value = 13

# Row 5 updates
$ws.Range("G5").Value = 2
$ws.Range("H5").Value = 3.5
$ws.Range("I5").Value = 3.7
$ws.Range("L5").Value = 4.5
$ws.Range("M5").Value = 1.08
$ws.Range("N5").Value = 8
$ws.Range("X5").Value = 8.5
$ws.Range("Y5").Value = 9
$ws.Range("AE5").Value = 17
$ws.Range("AO5").Value = 11
$ws.Range("BB5").Value = 251
